$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.257
$ws.Range("A6").Value = -22.291
$ws.Range("A7").Value = -20
$ws.Range("D7").Value = -8.154
$ws.Range("D12").Value = -7.302
$ws.Range("D15").Value = -8.196000000000002
$ws.Range("A16").Value = -21.779
$ws.Range("A20").Value = -20.14
$ws.Range("D20").Value = -7.714
$ws.Range("D21").Value = -7.988000000000001
$ws.Range("D22").Value = -7.970999999999999
$ws.Range("D23").Value = -7.997
$ws.Range("A28").Value = -22.111
$ws.Range("A29").Value = -21.343
$ws.Range("D29").Value = -7.572999999999999
$ws.Range("A32").Value = -21.648
$ws.Range("D34").Value = -7.903999999999999
$ws.Range("A40").Value = -19.896
$ws.Range("D42").Value = -7.952
$ws.Range("D43").Value = -7.853000000000002
$ws.Range("D44").Value = -7.831999999999999
$ws.Range("D45").Value = -7.531000000000001
$ws.Range("A46").Value = -21.786
$ws.Range("D46").Value = -8.374000000000001
$ws.Range("D50").Value = -8.355
$ws.Range("A51").Value = -22.162
$ws.Range("D51").Value = -8.301
$ws.Range("A52").Value = -22.23
$ws.Range("A57").Value = -22.263
$ws.Range("A59").Value = -22.354
$ws.Range("A62").Value = -22.085
$ws.Range("A66").Value = -21.532
$ws.Range("D66").Value = -7.557
$ws.Range("D67").Value = -7.204000000000001
$ws.Range("A73").Value = -20.213
$ws.Range("A74").Value = -21.244
$ws.Range("D79").Value = -7.6
$ws.Range("D84").Value = -8.300000000000001
$ws.Range("A92").Value = -21.64
$ws.Range("D92").Value = -6.654000000000001
$ws.Range("D97").Value = -8.388
$ws.Range("A100").Value = -22.217
